$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) contains values that look numeric (e.g. "0.9999",
# "239.27") as well as values that are not valid numbers (e.g. "30.308.53",
# containing two dots). In the source data ALL of these are plain text.
# Pre-format the Price column cells as Text so Excel stores the numeric-
# looking ones as literal strings too, instead of silently converting them
# to floating point numbers on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.308.53"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").Value = "1.888.68"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "239.27"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "0.4685"
$ws.Range("E7").Value = "  -1.82%  "

$ws.Range("D8").Value = "0.2872"
$ws.Range("E8").Value = "  +0.71%  "

$ws.Range("D9").Value = "0.06620"
$ws.Range("E9").Value = "  -1.01%  "

$ws.Range("D10").Value = "20.12"
$ws.Range("E10").Value = "  +6.91%  "

$ws.Range("D11").Value = "0.07783"
$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("D12").Value = "98.40"
$ws.Range("E12").Value = "  -3.08%  "

$ws.Range("D13").Value = "1.899.35"
$ws.Range("E13").Value = "  -0.86%  "

$ws.Range("D14").Value = "5.147"
$ws.Range("E14").Value = "  -1.89%  "

$ws.Range("D15").Value = "0.6846"
$ws.Range("E15").Value = "  +1.89%  "

$ws.Range("D16").Value = "285.35"
$ws.Range("E16").Value = "  +11.30%  "

$ws.Range("D17").Value = "30.319.40"
$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").Value = "0.9996"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "12.66"
$ws.Range("E19").Value = "  -0.30%  "

$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.141.24"
$ws.Range("E20").Value = "  -0.98%  "

$ws.Range("D21").Value = "5.398"
$ws.Range("E21").Value = "  -0.04%  "

$ws.Range("D22").Value = "0.000007336"
$ws.Range("E22").Value = "  -2.09%  "

$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").Value = "6.217"
$ws.Range("E24").Value = "  -1.47%  "

$ws.Range("D25").Value = "9.418"
$ws.Range("E25").Value = "  +0.70%  "

$ws.Range("D26").Value = "165.51"
$ws.Range("E26").Value = "  -1.98%  "

$ws.Range("D27").Value = "19.34"
$ws.Range("E27").Value = "  +1.46%  "

$ws.Range("D28").Value = "2.008"
$ws.Range("E28").Value = "  -3.13%  "

$ws.Range("D29").Value = "1.379"
$ws.Range("E29").Value = "  -0.38%  "

$ws.Range("D30").Value = "0.09774"
$ws.Range("E30").Value = "  -3.36%  "

$ws.Range("D31").Value = "4.484"
$ws.Range("E31").Value = "  -4.85%  "

$ws.Range("D32").Value = "1.489"
$ws.Range("E32").Value = "  -1.77%  "

$ws.Range("D33").Value = "4.186"
$ws.Range("E33").Value = "  -1.77%  "

$ws.Range("D34").Value = "0.04733"
$ws.Range("E34").Value = "  +0.12%  "

$ws.Range("D35").Value = "0.7159"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("D36").Value = "1.100"
$ws.Range("E36").Value = "  -0.92%  "

$ws.Range("D37").Value = "2.708"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").Value = "0.01887"
$ws.Range("E38").Value = "  -1.45%  "

$ws.Range("D39").Value = "6.660"
$ws.Range("E39").Value = "  +6.53%  "

$ws.Range("D40").Value = "2.537"
$ws.Range("E40").Value = "  -2.99%  "

$ws.Range("D41").Value = "73.05"
$ws.Range("E41").Value = "  -2.97%  "

$ws.Range("D42").Value = "1.986"
$ws.Range("E42").Value = "  +0.69%  "

$ws.Range("D43").Value = "0.8733"
$ws.Range("E43").Value = "  +1.23%  "

$ws.Range("D44").Value = "104.45"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4242"
$ws.Range("E45").Value = "  -0.28%  "

$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9994"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").Value = "988.51"
$ws.Range("E47").Value = "  +1.18%  "

$ws.Range("D48").Value = "7.299"
$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("D49").Value = "9.270"
$ws.Range("E49").Value = "  +4.66%  "

$ws.Range("D50").Value = "0.1167"
$ws.Range("E50").Value = "  -2.75%  "

$ws.Range("D51").Value = "34.24"
$ws.Range("E51").Value = "  -1.74%  "

# Restore the Price column cells to the default ("Normal") style so only
# their content changed and no stray number-format/style is left behind on
# cells -- the stored values remain text because they were already committed
# as text above; only the cosmetic number format is being reset here.
$ws.Range("D2:D51").Style = "Normal"
